$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 400
    4  = 800
    5  = 800
    6  = 800
    7  = 800
    8  = 800
    9  = 800
    12 = 300
    13 = 800
    14 = 800
    15 = 800
    16 = 300
    18 = 700
    19 = 800
    20 = 800
    21 = 800
    22 = 400
    23 = 800
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
